$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column H: sqrt(F0) for each data row of Table18 (B4:G12)
# Header cell H4 mirrors the formatting of the other header cells (row 4)
# but is left without a label, matching the source edit.
$ws.Cells.Item(4, 8).HorizontalAlignment = -4108

for ($r = 5; $r -le 12; $r++) {
    $cell = $ws.Cells.Item($r, 8)
    $cell.Formula = "=SQRT(Table18[[#This Row],[F0 in N]])"
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4108
}

# Move the view / selection the way the author left it: scrolled down one
# row, with the cursor resting just past the new column.
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
$ws.Range("I10").Select()
